$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2687
$ws1.Range("F3").Value = 583
$ws1.Range("F8").Value = 1250
$ws1.Range("F9").Value = 587
$ws1.Range("F14").Value = 5863
$ws1.Range("F15").Value = 99
$ws1.Range("F17").Value = 4299
$ws1.Range("F21").Value = 5002
$ws1.Range("F22").Value = 6448
$ws1.Range("F25").Value = 706
$ws1.Range("F31").Value = 1006
$ws1.Range("F32").Value = 1435
$ws1.Range("F37").Value = 1765
$ws1.Range("F39").Value = 1166
$ws1.Range("F40").Value = 1341
$ws1.Range("F44").Value = 139
$ws1.Range("F48").Value = 45

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 1218
$ws2.Range("F16").Value = 15

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 4055

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4055
$ws4.Range("F4").Value = 583
$ws4.Range("F7").Value = 1218
$ws4.Range("F12").Value = 1250
$ws4.Range("F14").Value = 587
$ws4.Range("F18").Value = 99
$ws4.Range("F20").Value = 4299
$ws4.Range("F21").Value = 5002
$ws4.Range("F24").Value = 706
$ws4.Range("F30").Value = 1435
$ws4.Range("F35").Value = 1765
$ws4.Range("F44").Value = 139
$ws4.Range("F48").Value = 45
